$d = $word.ActiveDocument

function Set-ParaText($paraIndex, $oldText, $newText) {
    $p = $d.Paragraphs($paraIndex)
    $full = $p.Range
    # Exclude the trailing paragraph mark from the range so only the
    # visible text is replaced.
    $r = $d.Range($full.Start, $full.End - 1)
    if ($r.Text -ne $oldText) {
        throw "Paragraph $paraIndex text mismatch. Expected '$oldText' but found '$($r.Text)'"
    }
    $r.Text = $newText
}

# Main page title
Set-ParaText 1 `
    "Play Finn's Golden Tavern for Free - Innovative Spiral Grid Gameplay" `
    "Play Finn's Golden Tavern for Free"

# "What we like" bullet points
Set-ParaText 42 `
    "Innovative gameplay with a spiral grid and wild symbol" `
    "Completely original and innovative gameplay"

Set-ParaText 43 `
    "Excellent graphics with a detailed tavern background" `
    "Excellent graphics with detailed background and symbols"

Set-ParaText 44 `
    "Original and engaging with free spins and multiplier features" `
    "Incredibly original and fantasy-themed"

Set-ParaText 45 `
    "Compelling writing style that informs and engages readers" `
    "Compelling writing style that keeps readers engaged"

# "What we don't like" bullet points
Set-ParaText 47 `
    "Hard to categorize and find identical games for comparison" `
    "Hard to categorize in one set genre"

Set-ParaText 48 `
    "May not appeal to traditionalists who prefer standard reels" `
    "Limited number of similar games available"

# Bold title repeated near the end of the document
Set-ParaText 49 `
    "Play Finn's Golden Tavern for Free - Innovative Spiral Grid Gameplay" `
    "Play Finn's Golden Tavern for Free"

# Italic summary paragraph
Set-ParaText 50 `
    "Read our review of Finn's Golden Tavern and play it for free. Innovative spiral grid gameplay with excellent graphics and original features. " `
    "Read our review of Finn's Golden Tavern to learn more about the game and play it for free."
